$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushes existing rows 4..243 down to 5..244),
# mirroring the author inserting a new EB wave entry (ZA7901 / EB 97.4)
# just under the most recent waves already listed at the top of the sheet.
$ws.Rows.Item(4).Insert()

# Fill in the new row's data.
$ws.Range("A4").Value = "ZA7901"
# Leading apostrophe forces text storage (matches the quote-prefixed "wave"
# values like 97.1, 97.2, ... already used in this column) instead of Excel
# auto-converting "97.4" to a number.
$ws.Range("B4").Value = "'97.4"
$ws.Range("C4").Value = "May-June 2022"
$ws.Range("D4").Value = "Fairness perceptions of the green transition"

# Restore the author's final selection/cursor position.
$ws.Range("D5").Select() | Out-Null
